# "exclusion to white milk in ccanz"
#
# The "Exclude" sheet's category filters are cleaned up / normalised:
#   "FLAVOURED MILK,07- JUICE" (and its inconsistent spacing variants)
#       -> "Flavoured Milk,Juice"
#   "WHITE MILK,Irrelevant,General"
#       -> "Irrelevant,General,White milk"
# and the workbook is left with the "Exclude" sheet active/selected
# (cursor resting on F19) instead of "Include".

$wb = $excel.ActiveWorkbook
$exclude = $wb.Worksheets.Item("Exclude")
$include = $wb.Worksheets.Item("Include")

# Numerator / Denominator rows that reference the flavoured-milk & juice
# category (previously inconsistently capitalised / spaced).
$exclude.Range("D2").Value2 = "Flavoured Milk,Juice"
$exclude.Range("D4").Value2 = "Flavoured Milk,Juice"
$exclude.Range("D7").Value2 = "Flavoured Milk,Juice"
$exclude.Range("D9").Value2 = "Flavoured Milk,Juice"

# Numerator / Denominator rows that reference the white-milk category
# (reordered + "White milk" de-capitalised).
$exclude.Range("D3").Value2 = "Irrelevant,General,White milk"
$exclude.Range("D5").Value2 = "Irrelevant,General,White milk"
$exclude.Range("D6").Value2 = "Irrelevant,General,White milk"
$exclude.Range("D8").Value2 = "Irrelevant,General,White milk"

# "Include" sheet content is untouched.

# Make "Exclude" the active tab again and leave the selection where the
# author's cursor ended up.
$exclude.Activate()
$exclude.Range("F19").Select()
